$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update figures after 85% data collected (row 2 and row 3 values)
$ws.Range("B2").Value = 0.480645161290323
$ws.Range("C2").Value = 0.536327608982827
$ws.Range("D2").Value = 0.49171270718232
$ws.Range("E2").Value = 0.512422360248447
$ws.Range("F2").Value = 0.523148148148148

$ws.Range("B3").Value = 0.617741935483871
$ws.Range("C3").Value = 0.598414795244386
$ws.Range("D3").Value = 0.644567219152855
$ws.Range("E3").Value = 0.566770186335404
$ws.Range("F3").Value = 0.585905349794239
